$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell J1 = "Q8", copying formatting from I1 (bold, centered, bordered header style)
$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update data values (regenerated series) and extend rows 4 & 8 with additional columns
$ws.Range("B2").Value = -1.287171013491868
$ws.Range("C2").Value = 0.2827939217689028
$ws.Range("D2").Value = 2.496697106700115
$ws.Range("E2").Value = 1.910099417287288
$ws.Range("F2").Value = 0.2999931319823368
$ws.Range("G2").Value = -0.1890695212528862
$ws.Range("H2").Value = 1.812607396685326
$ws.Range("B3").Value = 0.5817192572292627
$ws.Range("C3").Value = 2.795622442160475
$ws.Range("D3").Value = 2.209024752747648
$ws.Range("E3").Value = 0.5989184674426966
$ws.Range("F3").Value = 0.1098558142074737
$ws.Range("G3").Value = 2.111532732145686
$ws.Range("B4").Value = 2.860352663279089
$ws.Range("C4").Value = 2.273754973866262
$ws.Range("D4").Value = 0.6636486885613111
$ws.Range("E4").Value = 0.1745860353260881
$ws.Range("F4").Value = 2.1762629532643
$ws.Range("G4").Value = 0.9355192051181277
$ws.Range("H4").Value = -0.004342528784518129
$ws.Range("I4").Value = 1.104249791356595
$ws.Range("J4").Value = 0.5344944032044353
$ws.Range("B5").Value = 3.993739063922089
$ws.Range("C5").Value = 2.383632778617138
$ws.Range("D5").Value = 1.894570125381915
$ws.Range("E5").Value = 3.896247043320127
$ws.Range("F5").Value = 2.655503295173955
$ws.Range("G5").Value = 1.715641561271309
$ws.Range("H5").Value = 2.824233881412422
$ws.Range("I5").Value = 2.254478493260262
$ws.Range("B6").Value = 1.327939736723062
$ws.Range("C6").Value = 0.838877083487839
$ws.Range("D6").Value = 2.840554001426051
$ws.Range("E6").Value = 1.599810253279879
$ws.Range("F6").Value = 0.6599485193772328
$ws.Range("G6").Value = 1.768540839518346
$ws.Range("H6").Value = 1.198785451366186
$ws.Range("B7").Value = 0.596462224740588
$ws.Range("C7").Value = 2.5981391426788
$ws.Range("D7").Value = 1.357395394532628
$ws.Range("E7").Value = 0.4175336606299818
$ws.Range("F7").Value = 1.526125980771095
$ws.Range("G7").Value = 0.9563705926189352
$ws.Range("B8").Value = 2.147895666590351
$ws.Range("C8").Value = 0.9071519184441783
$ws.Range("D8").Value = -0.03270981545846752
$ws.Range("E8").Value = 1.075882504682646
$ws.Range("F8").Value = 0.5061271165304859
$ws.Range("G8").Value = 0.8472646462226479
$ws.Range("H8").Value = 0.7681077429517771
$ws.Range("I8").Value = 0.5692784093517318
$ws.Range("B9").Value = 0.2723519233738259
$ws.Range("C9").Value = -0.66750981052882
$ws.Range("D9").Value = 0.4410825096122935
$ws.Range("E9").Value = -0.1286728785398665
$ws.Range("F9").Value = 0.2124646511522954
$ws.Range("G9").Value = 0.1333077478814246
$ws.Range("H9").Value = -0.0655215857186206
$ws.Range("B10").Value = -0.6913225886447478
$ws.Range("C10").Value = 0.4172697314963657
$ws.Range("D10").Value = -0.1524856566557943
$ws.Range("E10").Value = 0.1886518730363676
$ws.Range("F10").Value = 0.1094949697654968
$ws.Range("G10").Value = -0.0893343638345484
$ws.Range("B11").Value = 0.2489548475442122
$ws.Range("C11").Value = -0.3208005406079479
$ws.Range("D11").Value = 0.02033698908421404
$ws.Range("E11").Value = -0.05881991418665673
$ws.Range("F11").Value = -0.257649247786702
$ws.Range("B12").Value = -0.4452547693108855
$ws.Range("C12").Value = -0.1041172396187236
$ws.Range("D12").Value = -0.1832741428895944
$ws.Range("E12").Value = -0.3821034764896396
$ws.Range("B13").Value = -0.1205552980991258
$ws.Range("C13").Value = -0.1997122013699966
$ws.Range("D13").Value = -0.3985415349700418
$ws.Range("B14").Value = -0.379135158058041
$ws.Range("C14").Value = -0.5779644916580863
$ws.Range("B15").Value = -0.3996235179026385

$ws.Range("A1").Select()
